# Update the "Correspond Handoff Datetime" (D5) and
# "Correspond Handback DateTime" (G5) values on the zh-cn and de-de
# report sheets to reflect the newly generated handback report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-02-16 10:27:52"
$wsZhCn.Range("G5").Value = "2016-02-16 10:28:46"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-02-16 10:28:06"
$wsDeDe.Range("G5").Value = "2016-02-16 10:29:14"
